# Inventario.xlsx edit
#
# Table1 currently has 3 columns (ID, Nombre, Empresa) and 4 data rows
# (2530/Laptop/CLK, 2531/Monitor/CLK, 2532/Mouse/CLK, 2533/xd/CLK).
#
# Target state: 4 columns (ID, Nombre, Cantidad, Fecha de compra) with a
# single data row: 3880 / Miguel / 6 / 23/12/2023.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# --- Table column changes -------------------------------------------------

# Rename the 3rd column header from "Empresa" to "Cantidad".
$tbl.ListColumns.Item(3).Range.Cells(1, 1).Value = "Cantidad"

# Add a 4th table column and give it its header: "Fecha de compra".
# ListColumns.Add() automatically expands the table's ref (A1:C1 -> A1:D1).
$newCol = $tbl.ListColumns.Add()
$newCol.Range.Cells(1, 1).Value = "Fecha de compra"

# --- Row changes -----------------------------------------------------------

# Drop the old rows 3-5, keeping only row 2 (which gets overwritten below).
$ws.Rows("3:5").Delete()

# Write the single remaining data row. "3880" and "6" look numeric, so a
# plain .Value assignment would store them as numbers; format the cells as
# text first (then clear that temporary format) so they are kept as text,
# matching the rest of the ID-like text data in the sheet.
$c = $ws.Cells.Item(2, 1)
$c.NumberFormat = "@"
$c.Value = "3880"
$c.ClearFormats()

$ws.Cells.Item(2, 2).Value = "Miguel"

$c = $ws.Cells.Item(2, 3)
$c.NumberFormat = "@"
$c.Value = "6"
$c.ClearFormats()

$ws.Cells.Item(2, 4).Value = "23/12/2023"
